# IST price update 2025-12-22 22:46
# Insert a new "latest snapshot" column before column B, shifting the
# existing price-history columns (B:AD) one column to the right (C:AE).
# The freshly-inserted column B gets the new timestamp header and, for
# every product row, the most recently known price (copied from the
# column that is now immediately to its right, i.e. the old column B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B; existing B:AD data moves to C:AE.
$ws.Columns("B:B").Insert()

# New timestamp header for the freshly inserted column.
$ws.Range("B1").Value = "2025-12-23 04:13"

# Determine the last used data row (SKU rows start at row 2).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Populate B2:B<lastRow> with the latest known price for each SKU, which
# is simply the value that now sits in column C (the column that used
# to be B before the insert).
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value2 = $ws.Cells.Item($r, 3).Value2
}
